$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 9 (pushes old rows 9.. down by one, and automatically
# extends the D8:D9 merged cell to D8:D10, and shifts D10:D11 -> D11:D12)
$ws.Rows("9:9").Insert()

# --- Row 8: fill in the previously-empty F/G/H cells for the "9" (tanpa
# epsilon) group now that its run has finished ---
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = "32m 59.36s"

# --- New row 9: epsilon = 1% sub-row inserted inside the merged D8:D10 cell ---

# D9 is the middle row of the 3-row vertical merge D8:D10, so it only needs
# a left/right border (top comes from D8's style, bottom from D10's style).
$d9 = $ws.Range("D9")
$d9.HorizontalAlignment = -4108   # xlCenter
$d9.VerticalAlignment = -4108     # xlCenter
$d9.Borders.Item(7).LineStyle = 1    # xlEdgeLeft
$d9.Borders.Item(10).LineStyle = 1   # xlEdgeRight

# E9 holds the epsilon percentage (1%), left aligned, no borders.
$e9 = $ws.Range("E9")
$e9.Value = 0.01
$e9.NumberFormat = "0%"
$e9.HorizontalAlignment = -4131   # xlLeft
$e9.VerticalAlignment = -4107     # xlBottom (default)
$e9.Borders.LineStyle = 0

# F9/G9/H9 carry the results for this epsilon row, boxed like the other
# data cells (same look as F4/G4/H4 etc, style index 2).
$ws.Range("F4").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Value = 70
$ws.Range("G4").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Value = 286
$ws.Range("H4").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("H9").Value = "19m 57.09s"

# --- Row 13 (previously row 12): E column is no longer used, drop it ---
$ws.Range("E13").Clear()

# --- Append a new blank row 19 (E:H only, same box style as row 18) ---
$ws.Range("E18:H18").Copy()
$ws.Range("E19:H19").PasteSpecial(-4122)

$excel.CutCopyMode = 0
